# Update the "TestData" sheet: rename the deal's negotiator and network.
#   D2: "LEGAL"  -> "Dan Harrison"
#   E2: "BRAVO"  -> "Universal HD"
# (Setting E2 before D2 keeps the shared-string insertion order aligned
#  with the target workbook: "Universal HD" then "Dan Harrison".)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

$ws.Range("E2").Value = "Universal HD"
$ws.Range("D2").Value = "Dan Harrison"

# Leave the selection on the last-edited cell, matching the authored workbook.
$ws.Range("D2").Select() | Out-Null
